$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Mann-Whitney test results (alpha_MW p-values and significance)
# Column A (Nutrient) stays the same; only columns B (alpha_MW) and C (significant) change.

$ws.Range("B2").Value = 0.0115848507956819
$ws.Range("C2").Value = "yes"

$ws.Range("B3").Value = 0.259489220400994
$ws.Range("C3").Value = "no"

$ws.Range("B4").Value = 0.825838735231449
$ws.Range("C4").Value = "no"

$ws.Range("B5").Value = 0.0000000024759263985344
$ws.Range("C5").Value = "yes"

$ws.Range("B6").Value = 0.000000409238204905671
$ws.Range("C6").Value = "yes"

$ws.Range("B7").Value = 0.000376151153428099
$ws.Range("C7").Value = "yes"

$ws.Range("B8").Value = 0.000000421170671648155
$ws.Range("C8").Value = "yes"

$ws.Range("B9").Value = 0.00000000653948129179547
$ws.Range("C9").Value = "yes"

$ws.Range("B10").Value = 0.00000000000000351026476366919
$ws.Range("C10").Value = "yes"

$ws.Range("B11").Value = 0.0339233039070592
$ws.Range("C11").Value = "yes"

$ws.Range("B12").Value = 0.0141541530448233
$ws.Range("C12").Value = "yes"

$ws.Range("B13").Value = 0.000000043633789126392
$ws.Range("C13").Value = "yes"

$ws.Range("B14").Value = 0.0000000684745185171894
$ws.Range("C14").Value = "yes"
